$d = $word.ActiveDocument

$d.Content.Find.Execute("Methods: shower() wash_hands() use_bathroom()", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Methods: tempature_control() fill() flush()", 2)
